# Adds a new "نواقص" item row into the DaySale report, re-numbers the rows
# that follow, replaces the last item with a new product, refreshes the
# running total and updates the "printed at" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $fmt = $cell.NumberFormat
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = $fmt
}

# Insert a new blank row above row 10 ("فرشه اسنان POWER GOLD كبار"); this
# shifts rows 10-13 down to 11-14. Then clone the (now shifted) row 11's
# formatting/merges back into the freshly inserted row 10 by copying it -
# this keeps every style/border/merge identical to the rest of the table.
$ws.Rows.Item(10).Insert()
$ws.Range("A11:Q11").Copy($ws.Range("A10:Q10"))
$ws.Rows.Item(10).RowHeight = 24.75

# --- Row 10: brand-new item --------------------------------------------
$ws.Cells.Item(10, 1).Value = 4
Set-TextValue 10 3  "بودره نلج اكياس"
Set-TextValue 10 8  "23:0"
Set-TextValue 10 12 "0"
Set-TextValue 10 14 "5.00"
Set-TextValue 10 16 "5.0000"
Set-TextValue 10 17 "1:0"

# --- Row 11: previous row 10 content, renumbered ------------------------
$ws.Cells.Item(11, 1).Value = 5

# --- Row 12: replaces the old "معجون سيجنال 25 مل" line with the new item
$ws.Cells.Item(12, 1).Value = 6
Set-TextValue 12 3  "معجون اسنان ديبرودنت 13 مل"
Set-TextValue 12 8  "5:0"
Set-TextValue 12 12 "0"
Set-TextValue 12 14 "45.00"
Set-TextValue 12 16 "45.0000"
Set-TextValue 12 17 "1:0"

# --- Row 13: running total ------------------------------------------------
$ws.Cells.Item(13, 16).Value = 112.67
$ws.Rows.Item(13).RowHeight = 24.75

# --- Row 14: footer (timestamp / page number / credit) -------------------
Set-TextValue 14 1 "Friday, 18 July, 2025 4:51 PM"
